$d = $word.ActiveDocument

# Namespace prefix used for every InsertXML payload below.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParaRangeXml([object]$range, [string]$bodyInner) {
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData><w:document ' + $wns + '><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Work from the bottom of the document upward so earlier edits don't shift
# the paragraph indices used by edits still to come.
# ---------------------------------------------------------------------------

# --- "Others" skills paragraph (last paragraph): drop the lastRenderedPageBreak
#     (it moves up to the "Software development" heading instead) and extend
#     the sentence with the new items, dropping the trailing period.
#     NB: this is the very last paragraph in the body, so InsertXML must
#     target its content only (end - 1, excluding the final pilcrow) -
#     replacing the whole paragraph range (incl. the body's mandatory
#     trailing mark) leaves a stray empty paragraph behind.
$pOthersList = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRng = $d.Range($pOthersList.Range.Start, $pOthersList.Range.End - 1)
$body = '<w:p><w:r><w:t>Illustration, photography, videography, graphic design</w:t></w:r>' + `
        '<w:r><w:t>, product design, CNC, metalworking</w:t></w:r></w:p>'
Set-ParaRangeXml $lastRng $body

# --- Software development skills list: split out " Node.js," as its own run.
$pSoftwareList = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("HTML, CSS, JavaScript")) { $pSoftwareList = $p }
}
$body = '<w:p><w:r><w:t>HTML, CSS, JavaScript, git, Tailwind CSS,</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> Node.js,</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> React, Next.js, Framer Motion, </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>SvelteKit</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve">, Svelte, </w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>Vite</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t>, Vue.js, WebGL, Three.js, Python</w:t></w:r></w:p>'
Set-ParaRangeXml $pSoftwareList.Range $body

# --- "Software development" heading: the rendered page break now falls here.
$pSoftwareHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Software development") { $pSoftwareHeading = $p }
}
$body = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Software development</w:t></w:r></w:p>'
Set-ParaRangeXml $pSoftwareHeading.Range $body

# --- Budgetty project heading: add the "(richardfxr.com/projects/budgetty)" link.
$pBudgetty = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Budgetty") { $pBudgetty = $p }
}
$body = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Budgetty</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
        '<w:r><w:t>richardfxr.com/projects/</w:t></w:r>' + `
        '<w:proofErr w:type="spellStart"/><w:r><w:t>budgetty</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
        '<w:r><w:t>)</w:t></w:r></w:p>'
Set-ParaRangeXml $pBudgetty.Range $body

# --- Weather Dial project heading: add the "(weatherdial.richardfxr.com)" link.
$pWeatherDial = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Weather Dial") { $pWeatherDial = $p }
}
$body = '<w:p><w:r><w:t>Weather Dial</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' + `
        '<w:r><w:t>weatherdial.richardfxr.com</w:t></w:r>' + `
        '<w:r><w:t>)</w:t></w:r></w:p>'
Set-ParaRangeXml $pWeatherDial.Range $body

# --- Education block: "Rhode Island School of Design" through
#     "Bachelor of Fine Arts, Industrial Design" expand into six paragraphs
#     covering RISD (with dates) and the newly added BUA entry.
$pRisd = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Rhode Island School of Design") { $pRisd = $p }
}
$pBfa = $pRisd.Next()
$eduRange = $d.Range($pRisd.Range.Start, $pBfa.Range.End)
$body = '<w:p><w:r><w:t>Rhode Island School of Design</w:t></w:r><w:r><w:t xml:space="preserve"> (RISD)</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>2020-2024</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>Bachelor of Fine Arts, Industrial Design</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>Boston University Academy (BUA)</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t>2016-2020</w:t></w:r></w:p>' + `
        '<w:p><w:r><w:t xml:space="preserve">Graduated </w:t></w:r><w:r><w:t>summa cum laude</w:t></w:r></w:p>'
Set-ParaRangeXml $eduRange $body

# --- Objective paragraph: call out "web development" skills specifically.
$pObjective = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Seeking an internship")) { $pObjective = $p }
}
$body = '<w:p><w:r><w:t xml:space="preserve">Seeking an internship in digital user interface design where I can apply my design and </w:t></w:r>' + `
        '<w:r><w:t>web development</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> skills while working with a professional team.</w:t></w:r></w:p>'
Set-ParaRangeXml $pObjective.Range $body

# --- Intro paragraph: "my résumé" instead of "Richard Fu's résumé", updated
#     date, and a simplified closing sentence pointing at the linked PDF.
$pIntro = $d.Paragraphs.Item(1)
$body = '<w:p><w:r><w:t xml:space="preserve">This is the web version of </w:t></w:r>' + `
        '<w:r><w:t>my</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve"> r</w:t></w:r>' + `
        '<w:r><w:t>ésumé</w:t></w:r>' + `
        '<w:r><w:t>. It was last updated on October 1</w:t></w:r>' + `
        '<w:r><w:t>5</w:t></w:r>' + `
        '<w:r><w:t xml:space="preserve">, 2022. </w:t></w:r>' + `
        '<w:r><w:t>Please use the PDF version linked below for downloading and printing.</w:t></w:r></w:p>'
Set-ParaRangeXml $pIntro.Range $body

Write-Host "Final paragraph count: $($d.Paragraphs.Count)"
